# recount | graphics fix
# Update the recomputed values in the data table (matlab service export).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = [double]"0.99192006418782341"
$ws.Range("B2").Value = [double]"309.04873793997916"
$ws.Range("C2").Value = [double]"358.78133278908865"
$ws.Range("D2").Value = [double]"27.841416924210222"
$ws.Range("E2").Value = [double]"32.321726780561526"
$ws.Range("F2").Value = [double]"3.9965666361292035E-2"
$ws.Range("G2").Value = [double]"0.81809763313609429"
$ws.Range("I2").Value = [double]"156752862.81520677"

$ws.Range("A3").Value = [double]"0.4157699305000736"
$ws.Range("B3").Value = [double]"0.6361433039288471"
$ws.Range("C3").Value = [double]"50.649184543720224"
$ws.Range("D3").Value = [double]"5.8158574597714008E-2"
$ws.Range("E3").Value = [double]"4.5899016244759911"
$ws.Range("F3").Value = [double]"4.9881197410415891E-2"
$ws.Range("I3").Value = [double]"271509454.42879027"

$ws.Range("A4").Value = [double]"0.99302983780383469"
$ws.Range("B4").Value = [double]"493.29602857262228"
$ws.Range("C4").Value = [double]"543.36442854871268"
$ws.Range("D4").Value = [double]"44.176281137638618"
$ws.Range("E4").Value = [double]"48.660122922003112"
$ws.Range("F4").Value = [double]"3.9668328896154847E-2"
$ws.Range("G4").Value = [double]"0.82415680473372743"
$ws.Range("I4").Value = [double]"195140933.00936475"

$ws.Range("A5").Value = [double]"0.41404917254056373"
$ws.Range("B5").Value = [double]"0.6639152696053543"
$ws.Range("C5").Value = [double]"50.563002234733545"
$ws.Range("D5").Value = [double]"6.0541303868019369E-2"
$ws.Range("E5").Value = [double]"4.5757462277811571"
$ws.Range("F5").Value = [double]"4.9787033783367851E-2"
$ws.Range("I5").Value = [double]"271510276.55927974"

$ws.Range("A6").Value = [double]"0.99567415815144966"
$ws.Range("B6").Value = [double]"369.42540183886723"
$ws.Range("C6").Value = [double]"439.19591220053155"
$ws.Range("D6").Value = [double]"29.01900123566146"
$ws.Range("E6").Value = [double]"34.500907330906742"
$ws.Range("F6").Value = [double]"2.8590665341237886E-2"
$ws.Range("G6").Value = [double]"0.59371893491124261"
$ws.Range("I6").Value = [double]"159815916.07801786"

$ws.Range("A7").Value = [double]"0.58209682934979323"
$ws.Range("B7").Value = [double]"4.9649386060883103"
$ws.Range("C7").Value = [double]"75.215080466535142"
$ws.Range("D7").Value = [double]"0.46776842509143535"
$ws.Range("E7").Value = [double]"7.0530153530681581"
$ws.Range("F7").Value = [double]"4.9717410754572852E-2"
$ws.Range("I7").Value = [double]"273069668.62597317"

$ws.Range("A8").Value = [double]"0.99556162762487732"
$ws.Range("B8").Value = [double]"566.64059494631329"
$ws.Range("C8").Value = [double]"636.60514534259335"
$ws.Range("D8").Value = [double]"44.535249379858271"
$ws.Range("E8").Value = [double]"50.035471042989293"
$ws.Range("F8").Value = [double]"2.8487615738319504E-2"
$ws.Range("G8").Value = [double]"0.6010650887573955"
$ws.Range("I8").Value = [double]"196284467.6907362"

$ws.Range("A9").Value = [double]"0.57946601858870472"
$ws.Range("B9").Value = [double]"4.8851188109151735"
$ws.Range("C9").Value = [double]"74.82581593753406"
$ws.Range("D9").Value = [double]"0.46037536729605755"
$ws.Range("E9").Value = [double]"7.0065957982973295"
$ws.Range("F9").Value = [double]"4.970788831558088E-2"
$ws.Range("I9").Value = [double]"273040941.80044872"

$ws.Range("A11").Value = [double]"0.71716407280149264"
$ws.Range("B11").Value = [double]"23.220220952013833"
$ws.Range("C11").Value = [double]"84.359534252437641"
$ws.Range("D11").Value = [double]"2.0494763040423134"
$ws.Range("E11").Value = [double]"7.4457890313667434"
$ws.Range("F11").Value = [double]"4.6919995275545252E-2"
$ws.Range("I11").Value = [double]"181606327.00365141"

